# Adds a new "Ganancia" column (R) to the report sheet and appends four
# new prediction rows (96-99) that include the computed Ganancia value
# for each played number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: R1 = "Ganancia" -----------------------------------
$ws.Range("R1").Value = "Ganancia"

# Copy the header formatting (bold font, borders, centered alignment) from
# the existing header cells (e.g. A1) onto the new header cell so it matches
# the look of the rest of the header row.
$ws.Range("A1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New data rows 96-99 ----------------------------------------------------

# Row 96
$ws.Range("A96").Value = "2024-01-17 02:36:41"
$ws.Range("B96").Value = 8
$ws.Range("C96").Value = 7
$ws.Range("D96").Value = 1
$ws.Range("E96").Value = 2
$ws.Range("F96").Value = 2
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0.001
$ws.Range("J96").Value = 0.05
$ws.Range("K96").Value = 0.003
$ws.Range("L96").Value = 100
$ws.Range("M96").Value = 500
$ws.Range("N96").Value = 10
$ws.Range("O96").Value = 5
$ws.Range("Q96").Value = "Data/bombay1.xlsx"
$ws.Range("R96").Value = 126000

# Row 97
$ws.Range("A97").Value = "2024-01-19 01:32:00"
$ws.Range("B97").Value = 55
$ws.Range("C97").Value = 40
$ws.Range("D97").Value = 2
$ws.Range("E97").Value = 15
$ws.Range("F97").Value = 11
$ws.Range("G97").Value = 12
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0.001
$ws.Range("J97").Value = 0.05
$ws.Range("K97").Value = 0.003
$ws.Range("L97").Value = 100
$ws.Range("M97").Value = 500
$ws.Range("N97").Value = 10
$ws.Range("O97").Value = 5
$ws.Range("Q97").Value = "Data/bombay1.xlsx"
$ws.Range("R97").Value = -298000

# Row 98
$ws.Range("A98").Value = "2024-01-19 01:58:03"
$ws.Range("B98").Value = 35
$ws.Range("C98").Value = 28
$ws.Range("D98").Value = 1
$ws.Range("E98").Value = 6
$ws.Range("F98").Value = 5
$ws.Range("G98").Value = 7
$ws.Range("H98").Value = 9
$ws.Range("I98").Value = 0.001
$ws.Range("J98").Value = 0.05
$ws.Range("K98").Value = 0.003
$ws.Range("L98").Value = 100
$ws.Range("M98").Value = 500
$ws.Range("N98").Value = 10
$ws.Range("O98").Value = 5
$ws.Range("Q98").Value = "Data/bombay1.xlsx"
$ws.Range("R98").Value = 504000

# Row 99
$ws.Range("A99").Value = "2024-01-19 07:54:21"
$ws.Range("B99").Value = 29
$ws.Range("C99").Value = 25
$ws.Range("D99").Value = 5
$ws.Range("E99").Value = 8
$ws.Range("F99").Value = 6
$ws.Range("G99").Value = 6
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0.001
$ws.Range("J99").Value = 0.05
$ws.Range("K99").Value = 0.003
$ws.Range("L99").Value = 100
$ws.Range("M99").Value = 500
$ws.Range("N99").Value = 10
$ws.Range("O99").Value = 5
$ws.Range("Q99").Value = "Data/bombay1.xlsx"
$ws.Range("R99").Value = 573000
